$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Croatia" sheet -----------------------------------
# The new sheet mirrors the "Slovakia" sheet's layout/styles exactly
# (same column widths, merged cells, row styles, etc.), so clone it and
# place the copy after "Spain" (the current last tab).
$slovakia = $wb.Worksheets.Item("Slovakia")
$spain    = $wb.Worksheets.Item("Spain")
$slovakia.Copy($null, $spain)

$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

# Fill in the Croatia-specific data (new shared strings get appended to
# the shared-string table in the order they are first written: the NGC
# code first, then the market name).
$croatia.Range("B4").Value = "NGC-3139/T2488"
$croatia.Range("B2").Value = "Croatia Market"

# Croatia becomes the active tab, with D13 selected.
$croatia.Activate()
$croatia.Range("D13").Select()

# --- 2. Slovakia's old selection is no longer "sticky" ------------------
# Its sheetView now just shows a whole-sheet selection (matching the
# other non-active sheets such as Germany/Portugal).
$slovakia.Activate()
$slovakia.Range("A1:XFD1048576").Select()

# --- 3. Re-activate Croatia so it ends up the workbook's active tab -----
$croatia.Activate()
